# Auto-generated Excel COM-interop script applying market-price/profit
# updates to the per-job "Ultros_Profits" worksheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) as produced by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets("ALC")
$ws.Range("H2").Value = 287.5
$ws.Range("I2").Value = 287.5
$ws.Range("K2").Value = 287.5
$ws.Range("M2").Value = -174.5
$ws.Range("H28").Value = 1502.88
$ws.Range("I28").Value = 1439.1177
$ws.Range("K28").Value = 1439.1177
$ws.Range("M28").Value = -954.1177
$ws.Range("H43").Value = 2975
$ws.Range("J43").Value = 3200
$ws.Range("L43").Value = 3200
$ws.Range("N43").Value = -3338
$ws.Range("H70").Value = 4298.85
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4298.85
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 12896.55
$ws.Range("N70").Value = -13436.55
$ws.Range("H73").Value = 4298.85
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4298.85
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 12896.55
$ws.Range("N73").Value = -14768.55
$ws.Range("H76").Value = 4000
$ws.Range("I76").Value = 5120.2
$ws.Range("J76").Value = 2599.75
$ws.Range("K76").Value = 5120.2
$ws.Range("L76").Value = 2599.75
$ws.Range("M76").Value = -4805.2
$ws.Range("N76").Value = -3229.75
$ws.Range("H79").Value = 4000
$ws.Range("I79").Value = 5120.2
$ws.Range("J79").Value = 2599.75
$ws.Range("K79").Value = 5120.2
$ws.Range("L79").Value = 2599.75
$ws.Range("M79").Value = -4028.2
$ws.Range("N79").Value = -4783.75
$ws.Range("H92").Value = 913.72
$ws.Range("J92").Value = 1051.7778
$ws.Range("L92").Value = 1051.7778
$ws.Range("N92").Value = -3547.7778
$ws.Range("H113").Value = 6346.905
$ws.Range("I113").Value = 4142.778
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 4142.778
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -888.7780000000002
$ws.Range("N113").Value = -14508
$ws.Range("H127").Value = 6701.72
$ws.Range("I127").Value = 1140.1428
$ws.Range("K127").Value = 3420.4284
$ws.Range("M127").Value = 1539.5716
$ws.Range("H137").Value = 1698.4
$ws.Range("I137").Value = 1213.3158
$ws.Range("K137").Value = 3639.9474
$ws.Range("M137").Value = -1089.9474
$ws.Range("H141").Value = 3422.2727
$ws.Range("I141").Value = 3347.1904
$ws.Range("K141").Value = 10041.5712
$ws.Range("M141").Value = -4861.5712
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets("ARM")
$ws.Range("H2").Value = 13382
$ws.Range("I2").Value = 14293.479
$ws.Range("K2").Value = 14293.479
$ws.Range("M2").Value = -14180.479
$ws.Range("H45").Value = 3152
$ws.Range("I45").Value = 2162.5
$ws.Range("J45").Value = 3481.8333
$ws.Range("K45").Value = 2162.5
$ws.Range("L45").Value = 3481.8333
$ws.Range("M45").Value = -1785.5
$ws.Range("N45").Value = -4235.8333
$ws.Range("H61").Value = 3110.84
$ws.Range("I61").Value = 1379.5333
$ws.Range("J61").Value = 5707.8
$ws.Range("K61").Value = 1379.5333
$ws.Range("L61").Value = 5707.8
$ws.Range("M61").Value = -1167.5333
$ws.Range("N61").Value = -6131.8
$ws.Range("H74").Value = 1815.5938
$ws.Range("I74").Value = 1822.4445
$ws.Range("K74").Value = 1822.4445
$ws.Range("M74").Value = -948.4445000000001
$ws.Range("H77").Value = 1815.5938
$ws.Range("I77").Value = 1822.4445
$ws.Range("K77").Value = 9112.2225
$ws.Range("M77").Value = -4744.2225
$ws.Range("H116").Value = 13382
$ws.Range("I116").Value = 14293.479
$ws.Range("K116").Value = 14293.479
$ws.Range("M116").Value = -11999.479
$ws.Range("H118").Value = 122999
$ws.Range("J118").Value = 122999
$ws.Range("L118").Value = 122999
$ws.Range("N118").Value = -126313
$ws.Range("H132").Value = 4071.2
$ws.Range("I132").Value = 3584
$ws.Range("J132").Value = 4666.6665
$ws.Range("K132").Value = 10752
$ws.Range("L132").Value = 13999.9995
$ws.Range("M132").Value = -8222
$ws.Range("N132").Value = -19059.9995
$ws.Range("H136").Value = 3110.84
$ws.Range("I136").Value = 1379.5333
$ws.Range("J136").Value = 5707.8
$ws.Range("K136").Value = 4138.5999
$ws.Range("L136").Value = 17123.4
$ws.Range("M136").Value = -1588.5999
$ws.Range("N136").Value = -22223.4

# ---- BSM ----
$ws = $wb.Worksheets("BSM")
$ws.Range("H3").Value = 13382
$ws.Range("I3").Value = 14293.479
$ws.Range("K3").Value = 14293.479
$ws.Range("M3").Value = -14179.479
$ws.Range("H22").Value = 15261.375
$ws.Range("I22").Value = 17384.428
$ws.Range("K22").Value = 17384.428
$ws.Range("M22").Value = -17211.428

# ---- CRP ----
$ws = $wb.Worksheets("CRP")
$ws.Range("H22").Value = 799
$ws.Range("I22").Value = 799
$ws.Range("K22").Value = 799
$ws.Range("M22").Value = -449
$ws.Range("H58").Value = 2416.72
$ws.Range("I58").Value = 1993.7142
$ws.Range("J58").Value = 2955.0908
$ws.Range("K58").Value = 1993.7142
$ws.Range("L58").Value = 2955.0908
$ws.Range("M58").Value = -1790.7142
$ws.Range("N58").Value = -3361.0908
$ws.Range("H92").Value = 56533.668
$ws.Range("J92").Value = 56533.668
$ws.Range("L92").Value = 56533.668
$ws.Range("N92").Value = -61525.668
$ws.Range("H125").Value = 98746.75
$ws.Range("J125").Value = 98746.75
$ws.Range("L125").Value = 98746.75
$ws.Range("N125").Value = -103666.75
$ws.Range("H134").Value = 3199.5925
$ws.Range("J134").Value = 5811
$ws.Range("L134").Value = 17433
$ws.Range("N134").Value = -22503
$ws.Range("H136").Value = 2416.72
$ws.Range("I136").Value = 1993.7142
$ws.Range("J136").Value = 2955.0908
$ws.Range("K136").Value = 5981.142599999999
$ws.Range("L136").Value = 8865.2724
$ws.Range("M136").Value = -3431.142599999999
$ws.Range("N136").Value = -13965.2724

# ---- CUL ----
$ws = $wb.Worksheets("CUL")
$ws.Range("H137").Value = 1476.4445
$ws.Range("I137").Value = 1583.7142
$ws.Range("J137").Value = 1101
$ws.Range("K137").Value = 4751.142599999999
$ws.Range("L137").Value = 3303
$ws.Range("M137").Value = 348.8574000000008
$ws.Range("N137").Value = -13503

# ---- GSM ----
$ws = $wb.Worksheets("GSM")
$ws.Range("H70").Value = 380370.34
$ws.Range("I70").Value = 560555.5
$ws.Range("K70").Value = 560555.5
$ws.Range("M70").Value = -560285.5
$ws.Range("H73").Value = 380370.34
$ws.Range("I73").Value = 560555.5
$ws.Range("K73").Value = 560555.5
$ws.Range("M73").Value = -559619.5
$ws.Range("H107").Value = 580.94116
$ws.Range("J107").Value = 768.9
$ws.Range("L107").Value = 768.9
$ws.Range("N107").Value = -4608.9
$ws.Range("H113").Value = 11098.875
$ws.Range("I113").Value = 8194.5
$ws.Range("K113").Value = 8194.5
$ws.Range("M113").Value = -6024.5
$ws.Range("H132").Value = 2679.606
$ws.Range("I132").Value = 1456.8334
$ws.Range("K132").Value = 4370.5002
$ws.Range("M132").Value = -1840.5002

# ---- LTW ----
$ws = $wb.Worksheets("LTW")
$ws.Range("H55").Value = 1151.591
$ws.Range("I55").Value = 1349.1875
$ws.Range("K55").Value = 1349.1875
$ws.Range("M55").Value = -1176.1875
$ws.Range("H93").Value = 3386
$ws.Range("I93").Value = 3383.6924
$ws.Range("K93").Value = 3383.6924
$ws.Range("M93").Value = -2135.6924
$ws.Range("H118").Value = 43384.23
$ws.Range("J118").Value = 43384.23
$ws.Range("L118").Value = 43384.23
$ws.Range("N118").Value = -46698.23
$ws.Range("H122").Value = 4217.5
$ws.Range("I122").Value = 2490
$ws.Range("J122").Value = 5945
$ws.Range("K122").Value = 7470
$ws.Range("L122").Value = 17835
$ws.Range("M122").Value = -5020
$ws.Range("N122").Value = -22735
$ws.Range("H136").Value = 4849.759
$ws.Range("I136").Value = 3077.4443
$ws.Range("K136").Value = 9232.332900000001
$ws.Range("M136").Value = -6682.332900000001

# ---- WVR ----
$ws = $wb.Worksheets("WVR")
$ws.Range("H47").Value = 36247.5
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 36247.5
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 36247.5
$ws.Range("N47").Value = -37391.5
$ws.Range("H123").Value = 49888.816
$ws.Range("J123").Value = 49888.816
$ws.Range("L123").Value = 49888.816
$ws.Range("N123").Value = -59688.816
$ws.Range("H132").Value = 1823.4231
$ws.Range("I132").Value = 1162.0476
$ws.Range("J132").Value = 4601.2
$ws.Range("K132").Value = 3486.142800000001
$ws.Range("L132").Value = 13803.6
$ws.Range("M132").Value = -956.1428000000005
$ws.Range("N132").Value = -18863.6
$ws.Range("H136").Value = 2031.24
$ws.Range("I136").Value = 1114.4615
$ws.Range("K136").Value = 3343.3845
$ws.Range("M136").Value = -793.3844999999997
$ws.Range("M47").ClearContents()
